# Natmi following Dr Hou advice
#
# The NATMI LR-pair table for Tnfsf8-Tnfrsf8 is regenerated: a new
# "ECs -> ECs" sending/target-cluster row is introduced ahead of the
# existing "sCs -> ECs" row, and the statistics for the (now second)
# "sCs -> ECs" row are refreshed to reflect the updated pipeline run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster "ECs", Ligand "Tnfsf8", Receptor "Tnfrsf8", Target cluster "ECs"
$row2 = @(
    "ECs", "Tnfsf8", "Tnfrsf8", "ECs",
    1, 0.3333333333333333, 0.01729866666666667, 0.051896,
    0.0536668448805226, 0.0536668448805226,
    3, 1, 1.977575666666667, 5.932727,
    1, 1, 0.03420942226577778, 0.307884800392,
    0.0536668448805226, 0.0536668448805226
)

# Row 3: Sending cluster "sCs", Ligand "Tnfsf8", Receptor "Tnfrsf8", Target cluster "ECs"
$row3 = @(
    "sCs", "Tnfsf8", "Tnfrsf8", "ECs",
    3, 1, 0.3050356666666666, 0.915107,
    0.9463331551194774, 0.9463331551194774,
    3, 1, 1.977575666666667, 5.932727,
    1, 1, 0.6032311118654444, 5.429080006788999,
    0.9463331551194774, 0.9463331551194774
)

for ($col = 1; $col -le $row2.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $row2[$col - 1]
    $ws.Cells.Item(3, $col).Value = $row3[$col - 1]
}
